$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(15).Copy()
$ws.Rows.Item(14).PasteSpecial(-4122)
